$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source export got re-sorted upstream: the observation records that used to sit
# on row 3 and row 6 swapped places, and likewise the records on row 5 and row 7
# swapped places (row 4's record is untouched). Rather than moving whole rows
# around (which would risk Excel re-interpreting the untouched, shared text
# columns such as the date-like "2022-05-30" strings as real dates), only the
# cells whose values actually differ between the two records are updated here.

# --- Row 3 becomes the record that used to be on row 6 ---
$ws.Range("A3").Value = 111741120
$ws.Range("B3").Value = 56398
$ws.Range("E3").Value = 100109
$ws.Range("F3").Value = "Tretåig hackspett"
$ws.Range("G3").Value = "Picoides tridactylus"
$ws.Range("H3").Value = "(Linnaeus, 1758)"
$ws.Range("Q3").Value = 331468.5669229594
$ws.Range("R3").Value = 6627064.351006002
$ws.Range("J3").Value = ""
$ws.Range("AF3").Value = ""
$ws.Range("M3").Value = "färska spår"

# --- Row 5 becomes the record that used to be on row 7 ---
$ws.Range("A5").Value = 111741025
$ws.Range("B5").Value = 94134
$ws.Range("E5").Value = 53
$ws.Range("F5").Value = "Vedtrappmossa"
$ws.Range("G5").Value = "Crossocalyx hellerianus"
$ws.Range("H5").Value = "(Nees ex Lindenb.) Meyl."
$ws.Range("Q5").Value = 331437.2628167981
$ws.Range("R5").Value = 6627065.263253132
$ws.Range("L5").Value = ""

# --- Row 6 becomes the record that used to be on row 3 ---
$ws.Range("A6").Value = 111741014
$ws.Range("B6").Value = 94134
$ws.Range("E6").Value = 53
$ws.Range("F6").Value = "Vedtrappmossa"
$ws.Range("G6").Value = "Crossocalyx hellerianus"
$ws.Range("H6").Value = "(Nees ex Lindenb.) Meyl."
$ws.Range("Q6").Value = 331429.3527348472
$ws.Range("R6").Value = 6627058.050714097
$ws.Range("J6").Value = ""
$ws.Range("AF6").Value = ""
$ws.Range("M6").Value = ""

# --- Row 7 becomes the record that used to be on row 5 ---
$ws.Range("A7").Value = 111741082
$ws.Range("B7").Value = 77515
$ws.Range("E7").Value = 6425
$ws.Range("F7").Value = "Garnlav"
$ws.Range("G7").Value = "Alectoria sarmentosa"
$ws.Range("H7").Value = "(Ach.) Ach."
$ws.Range("Q7").Value = 331468.5669229594
$ws.Range("R7").Value = 6627064.351006002
$ws.Range("L7").Value = ""
